$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving the cells
# original style index (avoids leaving a stray NumberFormat behind).
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '25.628.04'
$ws.Range("E2").Value = '  -5.93%  '

Set-TextValue $ws.Range("D3") '1.806.87'
$ws.Range("E3").Value = '  -5.15%  '

$ws.Range("E4").Value = '  +0.11%  '

Set-TextValue $ws.Range("D5") '276.06'
$ws.Range("E5").Value = '  -9.87%  '

Set-TextValue $ws.Range("D6") '1.000'
$ws.Range("E6").Value = '  +0.12%  '

Set-TextValue $ws.Range("D7") '0.5071'
$ws.Range("E7").Value = '  -6.50%  '

Set-TextValue $ws.Range("D8") '0.3521'
$ws.Range("E8").Value = '  -7.57%  '

Set-TextValue $ws.Range("D9") '43.67'
$ws.Range("E9").Value = '  -5.09%  '

Set-TextValue $ws.Range("D10") '0.06629'
$ws.Range("E10").Value = '  -9.18%  '

Set-TextValue $ws.Range("D11") '20.03'
$ws.Range("E11").Value = '  -9.80%  '

Set-TextValue $ws.Range("D12") '0.8379'
$ws.Range("E12").Value = '  -7.37%  '

Set-TextValue $ws.Range("D13") '0.07767'
$ws.Range("E13").Value = '  -5.12%  '

Set-TextValue $ws.Range("D14") '1.796.77'
$ws.Range("E14").Value = '  +56.47%  '

Set-TextValue $ws.Range("D15") '5.084'
$ws.Range("E15").Value = '  -4.98%  '

Set-TextValue $ws.Range("D16") '87.63'
$ws.Range("E16").Value = '  -8.60%  '

$ws.Range("E17").Value = '  +0.07%  '

Set-TextValue $ws.Range("D18") '13.94'
$ws.Range("E18").Value = '  -6.30%  '

$ws.Range("E19").Value = '  +0.07%  '

Set-TextValue $ws.Range("D20") '0.000007959'
$ws.Range("E20").Value = '  -8.00%  '

Set-TextValue $ws.Range("D21") '25.684.74'
$ws.Range("E21").Value = '  -5.82%  '

Set-TextValue $ws.Range("D22") '4.726'
$ws.Range("E22").Value = '  -6.40%  '

Set-TextValue $ws.Range("D23") '10.03'
$ws.Range("E23").Value = '  -7.38%  '

Set-TextValue $ws.Range("D24") '6.059'
$ws.Range("E24").Value = '  -7.10%  '

Set-TextValue $ws.Range("D25") '142.62'
$ws.Range("E25").Value = '  -3.97%  '

Set-TextValue $ws.Range("D26") '2.109'
$ws.Range("E26").Value = '  -8.82%  '

Set-TextValue $ws.Range("D27") '1.656'
$ws.Range("E27").Value = '  -5.71%  '

Set-TextValue $ws.Range("D28") '16.94'
$ws.Range("E28").Value = '  -7.96%  '

Set-TextValue $ws.Range("D29") '108.24'
$ws.Range("E29").Value = '  -7.32%  '

Set-TextValue $ws.Range("D30") '4.325'
$ws.Range("E30").Value = '  -11.12%  '

$ws.Range("E31").Value = '  -10.34%  '

Set-TextValue $ws.Range("D32") '0.08793'
$ws.Range("E32").Value = '  -4.73%  '

Set-TextValue $ws.Range("D33") '0.04795'
$ws.Range("E33").Value = '  -5.72%  '

Set-TextValue $ws.Range("D34") '0.7254'
$ws.Range("E34").Value = '  -12.81%  '

Set-TextValue $ws.Range("D35") '1.127'
$ws.Range("E35").Value = '  -7.88%  '

Set-TextValue $ws.Range("D36") '2.857'
$ws.Range("E36").Value = '  -4.88%  '

Set-TextValue $ws.Range("D37") '1.0000'
$ws.Range("E37").Value = '  +0.17%  '

Set-TextValue $ws.Range("D38") '3.033'
$ws.Range("E38").Value = '  -8.67%  '

Set-TextValue $ws.Range("D39") '0.01862'
$ws.Range("E39").Value = '  -6.99%  '

Set-TextValue $ws.Range("D40") '0.5181'
$ws.Range("E40").Value = '  -12.86%  '

Set-TextValue $ws.Range("D41") '2.272'
$ws.Range("E41").Value = '  -15.79%  '

Set-TextValue $ws.Range("D42") '0.9623'
$ws.Range("E42").Value = '  -10.95%  '

Set-TextValue $ws.Range("D43") '114.79'
$ws.Range("E43").Value = '  -1.54%  '

Set-TextValue $ws.Range("D44") '6.179'
$ws.Range("E44").Value = '  -7.51%  '

Set-TextValue $ws.Range("D45") '8.032'
$ws.Range("E45").Value = '  -13.59%  '

$ws.Range("E46").Value = '  +0.14%  '

Set-TextValue $ws.Range("D47") '0.4569'
$ws.Range("E47").Value = '  -10.90%  '

Set-TextValue $ws.Range("D48") '0.1382'
$ws.Range("E48").Value = '  -9.72%  '

Set-TextValue $ws.Range("D49") '9.238'
$ws.Range("E49").Value = '  -9.39%  '

Set-TextValue $ws.Range("D50") '35.88'
$ws.Range("E50").Value = '  -6.41%  '

Set-TextValue $ws.Range("D51") '1.494'
$ws.Range("E51").Value = '  -9.13%  '
